$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Cells.Item(46, 8).Value = 1984.5
$ws.Cells.Item(46, 9).Value = 1984.5
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 5953.5
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -5834.5
$ws.Cells.Item(46, 14).ClearContents()
$ws.Cells.Item(60, 8).Value = 1984.5
$ws.Cells.Item(60, 9).Value = 1984.5
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 5953.5
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = -5469.5
$ws.Cells.Item(60, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 6437.8
$ws.Cells.Item(74, 9).Value = 5729.6665
$ws.Cells.Item(74, 10).Value = 7500
$ws.Cells.Item(74, 11).Value = 5729.6665
$ws.Cells.Item(74, 12).Value = 7500
$ws.Cells.Item(74, 13).Value = -4793.6665
$ws.Cells.Item(74, 14).Value = -9372
$ws.Cells.Item(77, 8).Value = 6437.8
$ws.Cells.Item(77, 9).Value = 5729.6665
$ws.Cells.Item(77, 10).Value = 7500
$ws.Cells.Item(77, 11).Value = 28648.3325
$ws.Cells.Item(77, 12).Value = 37500
$ws.Cells.Item(77, 13).Value = -23968.3325
$ws.Cells.Item(77, 14).Value = -46860
$ws.Cells.Item(100, 8).Value = 3081.348
$ws.Cells.Item(100, 10).Value = 3514.7778
$ws.Cells.Item(100, 12).Value = 3514.7778
$ws.Cells.Item(100, 14).Value = -4596.7778
$ws.Cells.Item(107, 8).Value = 2002
$ws.Cells.Item(107, 9).Value = 2002
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 2002
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -82
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 1770.0714
$ws.Cells.Item(135, 9).Value = 1771.6154
$ws.Cells.Item(135, 11).Value = 15944.5386
$ws.Cells.Item(135, 13).Value = -13409.5386
$ws.Cells.Item(137, 8).Value = 13094.167
$ws.Cells.Item(137, 9).Value = 1273.625
$ws.Cells.Item(137, 10).Value = 22550.6
$ws.Cells.Item(137, 11).Value = 3820.875
$ws.Cells.Item(137, 12).Value = 67651.79999999999
$ws.Cells.Item(137, 13).Value = -1270.875
$ws.Cells.Item(137, 14).Value = -72751.79999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Cells.Item(32, 8).Value = 3241.7346
$ws.Cells.Item(32, 9).Value = 2302.45
$ws.Cells.Item(32, 10).Value = 7416.3335
$ws.Cells.Item(32, 11).Value = 2302.45
$ws.Cells.Item(32, 12).Value = 7416.3335
$ws.Cells.Item(32, 13).Value = -2015.45
$ws.Cells.Item(32, 14).Value = -7990.3335
$ws.Cells.Item(61, 8).Value = 72774.85000000001
$ws.Cells.Item(61, 9).Value = 2439.4285
$ws.Cells.Item(61, 10).Value = 195861.83
$ws.Cells.Item(61, 11).Value = 2439.4285
$ws.Cells.Item(61, 12).Value = 195861.83
$ws.Cells.Item(61, 13).Value = -2227.4285
$ws.Cells.Item(61, 14).Value = -196285.83
$ws.Cells.Item(74, 8).Value = 13266.714
$ws.Cells.Item(74, 9).Value = 1756.0741
$ws.Cells.Item(74, 10).Value = 52115.125
$ws.Cells.Item(74, 11).Value = 1756.0741
$ws.Cells.Item(74, 12).Value = 52115.125
$ws.Cells.Item(74, 13).Value = -882.0741
$ws.Cells.Item(74, 14).Value = -53863.125
$ws.Cells.Item(76, 8).Value = 150000
$ws.Cells.Item(76, 10).Value = 150000
$ws.Cells.Item(76, 12).Value = 150000
$ws.Cells.Item(76, 14).Value = -150676
$ws.Cells.Item(77, 8).Value = 13266.714
$ws.Cells.Item(77, 9).Value = 1756.0741
$ws.Cells.Item(77, 10).Value = 52115.125
$ws.Cells.Item(77, 11).Value = 8780.370500000001
$ws.Cells.Item(77, 12).Value = 260575.625
$ws.Cells.Item(77, 13).Value = -4412.370500000001
$ws.Cells.Item(77, 14).Value = -269311.625
$ws.Cells.Item(79, 8).Value = 150000
$ws.Cells.Item(79, 10).Value = 150000
$ws.Cells.Item(79, 12).Value = 150000
$ws.Cells.Item(79, 14).Value = -152340
$ws.Cells.Item(132, 8).Value = 4183710.8
$ws.Cells.Item(132, 9).Value = 5831.6
$ws.Cells.Item(132, 11).Value = 17494.8
$ws.Cells.Item(132, 13).Value = -14964.8
$ws.Cells.Item(136, 8).Value = 72774.85000000001
$ws.Cells.Item(136, 9).Value = 2439.4285
$ws.Cells.Item(136, 10).Value = 195861.83
$ws.Cells.Item(136, 11).Value = 7318.2855
$ws.Cells.Item(136, 12).Value = 587585.49
$ws.Cells.Item(136, 13).Value = -4768.2855
$ws.Cells.Item(136, 14).Value = -592685.49

# --- Sheet: BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 13).ClearContents()
$ws.Cells.Item(57, 8).Value = 89997.5
$ws.Cells.Item(57, 10).Value = 89997.5
$ws.Cells.Item(57, 12).Value = 89997.5
$ws.Cells.Item(57, 14).Value = -91437.5
$ws.Cells.Item(105, 8).Value = 55557000
$ws.Cells.Item(105, 9).Value = 66668050
$ws.Cells.Item(105, 11).Value = 66668050
$ws.Cells.Item(105, 13).Value = -66666303
$ws.Cells.Item(134, 8).Value = 39927.184
$ws.Cells.Item(134, 9).Value = 36050.613
$ws.Cells.Item(134, 11).Value = 108151.839
$ws.Cells.Item(134, 13).Value = -105616.839
$ws.Cells.Item(136, 8).Value = 89997.5
$ws.Cells.Item(136, 10).Value = 89997.5
$ws.Cells.Item(136, 12).Value = 89997.5
$ws.Cells.Item(136, 14).Value = -100197.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Cells.Item(6, 8).Value = 8311100
$ws.Cells.Item(6, 9).Value = 11635040
$ws.Cells.Item(6, 11).Value = 11635040
$ws.Cells.Item(6, 13).Value = -11634927
$ws.Cells.Item(31, 8).Value = 17742.303
$ws.Cells.Item(31, 9).Value = 8467.799999999999
$ws.Cells.Item(31, 11).Value = 8467.799999999999
$ws.Cells.Item(31, 13).Value = -8172.799999999999
$ws.Cells.Item(34, 8).Value = 17742.303
$ws.Cells.Item(34, 9).Value = 8467.799999999999
$ws.Cells.Item(34, 11).Value = 8467.799999999999
$ws.Cells.Item(34, 13).Value = -8265.799999999999
$ws.Cells.Item(107, 8).Value = 717.5106
$ws.Cells.Item(107, 9).Value = 1271.2941
$ws.Cells.Item(107, 10).Value = 403.7
$ws.Cells.Item(107, 11).Value = 1271.2941
$ws.Cells.Item(107, 12).Value = 403.7
$ws.Cells.Item(107, 13).Value = 648.7058999999999
$ws.Cells.Item(107, 14).Value = -4243.7
$ws.Cells.Item(119, 8).Value = 44750
$ws.Cells.Item(119, 10).Value = 44750
$ws.Cells.Item(119, 12).Value = 44750
$ws.Cells.Item(119, 14).Value = -54426
$ws.Cells.Item(132, 8).Value = 1725.4166
$ws.Cells.Item(132, 9).Value = 1873.8
$ws.Cells.Item(132, 11).Value = 5621.4
$ws.Cells.Item(132, 13).Value = -3091.4
$ws.Cells.Item(134, 8).Value = 24395404
$ws.Cells.Item(134, 9).Value = 1789.9333
$ws.Cells.Item(134, 11).Value = 5369.7999
$ws.Cells.Item(134, 13).Value = -2834.7999

# --- Sheet: CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Cells.Item(104, 8).Value = 6174839.5
$ws.Cells.Item(104, 10).Value = 12345679
$ws.Cells.Item(104, 12).Value = 37037037
$ws.Cells.Item(104, 14).Value = -37042279

# --- Sheet: GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Cells.Item(2, 8).Value = 140.04347
$ws.Cells.Item(2, 9).Value = 84.05882
$ws.Cells.Item(2, 10).Value = 298.66666
$ws.Cells.Item(2, 11).Value = 84.05882
$ws.Cells.Item(2, 12).Value = 298.66666
$ws.Cells.Item(2, 13).Value = 28.94118
$ws.Cells.Item(2, 14).Value = -524.66666
$ws.Cells.Item(52, 8).Value = 972222340
$ws.Cells.Item(52, 10).Value = 1000000000
$ws.Cells.Item(52, 12).Value = 1000000000
$ws.Cells.Item(52, 14).Value = -1000000518
$ws.Cells.Item(74, 8).Value = 68500
$ws.Cells.Item(74, 10).Value = 68500
$ws.Cells.Item(74, 12).Value = 68500
$ws.Cells.Item(74, 14).Value = -70372
$ws.Cells.Item(77, 8).Value = 68500
$ws.Cells.Item(77, 10).Value = 68500
$ws.Cells.Item(77, 12).Value = 205500
$ws.Cells.Item(77, 14).Value = -214860
$ws.Cells.Item(113, 8).Value = 5007
$ws.Cells.Item(113, 9).Value = 4598.3335
$ws.Cells.Item(113, 10).Value = 5620
$ws.Cells.Item(113, 11).Value = 4598.3335
$ws.Cells.Item(113, 12).Value = 5620
$ws.Cells.Item(113, 13).Value = -2428.3335
$ws.Cells.Item(113, 14).Value = -9960
$ws.Cells.Item(122, 8).Value = 1360437.4
$ws.Cells.Item(122, 9).Value = 2265189
$ws.Cells.Item(122, 10).Value = 3309.7
$ws.Cells.Item(122, 11).Value = 6795567
$ws.Cells.Item(122, 12).Value = 9929.099999999999
$ws.Cells.Item(122, 13).Value = -6793117
$ws.Cells.Item(122, 14).Value = -14829.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Cells.Item(7, 8).Value = 3675522.8
$ws.Cells.Item(7, 9).Value = 5305533
$ws.Cells.Item(7, 11).Value = 5305533
$ws.Cells.Item(7, 13).Value = -5305421
$ws.Cells.Item(16, 8).Value = 142861550
$ws.Cells.Item(16, 9).Value = 142861550
$ws.Cells.Item(16, 11).Value = 142861550
$ws.Cells.Item(16, 13).Value = -142861380
$ws.Cells.Item(100, 8).Value = 5698.5454
$ws.Cells.Item(100, 9).Value = 6300.6665
$ws.Cells.Item(100, 11).Value = 6300.6665
$ws.Cells.Item(100, 13).Value = -5759.6665
$ws.Cells.Item(126, 8).Value = 3675522.8
$ws.Cells.Item(126, 9).Value = 5305533
$ws.Cells.Item(126, 11).Value = 15916599
$ws.Cells.Item(126, 13).Value = -15914129
$ws.Cells.Item(132, 8).Value = 2682930.5
$ws.Cells.Item(132, 10).Value = 6703210
$ws.Cells.Item(132, 12).Value = 20109630
$ws.Cells.Item(132, 14).Value = -20114690

# --- Sheet: WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Cells.Item(132, 8).Value = 20321.818
$ws.Cells.Item(132, 9).Value = 9530.583000000001
$ws.Cells.Item(132, 11).Value = 28591.749
$ws.Cells.Item(132, 13).Value = -26061.749
$ws.Cells.Item(136, 8).Value = 18369.393
$ws.Cells.Item(136, 10).Value = 47009.777
$ws.Cells.Item(136, 12).Value = 141029.331
$ws.Cells.Item(136, 14).Value = -146129.331
